$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 347) holds the "Förändrad" date, stored as the
# Excel date serial 45186 (2023-09-17). Update it to serial 45188
# (2023-09-19) for every data row, matching the original cell formatting.
$ws.Range("C2:C347").Value = 45188
